$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 766341
$ws.Range("D2").Value = 155968
$ws.Range("E2").Value = 1429248811

$ws.Range("C13").Value = 187861
$ws.Range("D13").Value = 33264
$ws.Range("E13").Value = 1168232002

$ws.Range("C19").Value = 27514
$ws.Range("D19").Value = 4210
$ws.Range("E19").Value = 132331474

$ws.Range("C21").Value = 175242
$ws.Range("D21").Value = 38057
$ws.Range("E21").Value = 316822468

$ws.Range("C41").Value = 126946
$ws.Range("D41").Value = 24257
$ws.Range("E41").Value = 662711335

$ws.Range("C57").Value = 31597
$ws.Range("D57").Value = 6018
$ws.Range("E57").Value = 162625330

$ws.Range("C81").Value = 88357
$ws.Range("D81").Value = 16598
$ws.Range("E81").Value = 499681283

$ws.Range("C88").Value = 71278
$ws.Range("D88").Value = 12436
$ws.Range("E88").Value = 110315721

$ws.Range("C121").Value = 1306396
$ws.Range("D121").Value = 220388
$ws.Range("E121").Value = 2275480279

$ws.Range("C129").Value = 633775
$ws.Range("D129").Value = 104970
$ws.Range("E129").Value = 3435036140

$ws.Range("C132").Value = 586022
$ws.Range("D132").Value = 90786
$ws.Range("E132").Value = 3472775626

$ws.Range("C151").Value = 39935
$ws.Range("D151").Value = 7155
$ws.Range("E151").Value = 60390275

$ws.Range("C156").Value = 12412
$ws.Range("D156").Value = 2144
$ws.Range("E156").Value = 40642757

$ws.Range("C171").Value = 95828
$ws.Range("D171").Value = 18309
$ws.Range("E171").Value = 490704985

$ws.Range("C178").Value = 515888
$ws.Range("D178").Value = 115380
$ws.Range("E178").Value = 891213559

$ws.Range("C186").Value = 236840
$ws.Range("D186").Value = 46319
$ws.Range("E186").Value = 1190035149

$ws.Range("C237").Value = 283328
$ws.Range("D237").Value = 49797
$ws.Range("E237").Value = 1438555325

$ws.Range("C240").Value = 205932
$ws.Range("D240").Value = 33984
$ws.Range("E240").Value = 1069916235
